# Updated symbol list for rows 2-51: refreshed price (D), 1h volume % (E),
# and hour (G) columns, matching the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Addr, $NewValue) {
    # Force the cell to Text format before assigning so the numeric-looking
    # literal (price / "12.34%" / hour) is kept verbatim instead of being
    # reinterpreted by Excel as a Number/Percentage/Date.
    $range = $Sheet.Range($Addr)
    $range.NumberFormat = "@"
    $range.Value = $NewValue
}

$edits = @(
    @("D2", "328.02"),
    @("E2", "5.90%"),
    @("G2", "15"),
    @("D3", "39.85"),
    @("E3", "7.06%"),
    @("G3", "15"),
    @("D4", "5.645"),
    @("E4", "10.01%"),
    @("G4", "15"),
    @("D5", "0.08089"),
    @("E5", "3.27%"),
    @("G5", "15"),
    @("D6", "4.547"),
    @("E6", "3.30%"),
    @("G6", "15"),
    @("D7", "8.678"),
    @("E7", "4.83%"),
    @("G7", "15"),
    @("D8", "1.954"),
    @("E8", "3.91%"),
    @("G8", "15"),
    @("G9", "15"),
    @("D10", "0.9491"),
    @("E10", "2.60%"),
    @("G10", "15"),
    @("D11", "0.1281"),
    @("E11", "8.95%"),
    @("G11", "15"),
    @("D12", "0.1987"),
    @("E12", "5.60%"),
    @("G12", "15"),
    @("D13", "0.09201"),
    @("E13", "4.09%"),
    @("G13", "15"),
    @("D14", "0.03560"),
    @("E14", "7.56%"),
    @("G14", "15"),
    @("D15", "0.09601"),
    @("E15", "-0.05%"),
    @("G15", "15"),
    @("D16", "0.001323"),
    @("E16", "-4.69%"),
    @("G16", "15"),
    @("D17", "0.006117"),
    @("E17", "-1.30%"),
    @("G17", "15"),
    @("D18", "3.373"),
    @("E18", "-0.59%"),
    @("G18", "15"),
    @("D19", "0.3505"),
    @("E19", "1.35%"),
    @("G19", "15"),
    @("D20", "7.492"),
    @("E20", "17.30%"),
    @("G20", "15"),
    @("E21", "8.40%"),
    @("G21", "15"),
    @("D22", "0.2507"),
    @("E22", "4.12%"),
    @("G22", "15"),
    @("D23", "0.04437"),
    @("E23", "2.18%"),
    @("G23", "15"),
    @("D24", "0.001256"),
    @("E24", "4.66%"),
    @("G24", "15"),
    @("D25", "0.004315"),
    @("E25", "0.70%"),
    @("G25", "15"),
    @("D26", "0.0001193"),
    @("E26", "-14.82%"),
    @("G26", "15"),
    @("D27", "0.0004001"),
    @("E27", "37.91%"),
    @("G27", "15"),
    @("G28", "15"),
    @("G29", "15"),
    @("G30", "15"),
    @("G31", "15"),
    @("G32", "15"),
    @("G33", "15"),
    @("G34", "15"),
    @("G35", "15"),
    @("G36", "15"),
    @("G37", "15"),
    @("G38", "15"),
    @("D39", "0.02514"),
    @("E39", "16.43%"),
    @("G39", "15"),
    @("D40", "0.05205"),
    @("E40", "3.90%"),
    @("G40", "15"),
    @("D41", "0.007810"),
    @("E41", "2.92%"),
    @("G41", "15"),
    @("D42", "0.1433"),
    @("E42", "5.62%"),
    @("G42", "15"),
    @("D43", "0.009074"),
    @("E43", "6.94%"),
    @("G43", "15"),
    @("D44", "0.002195"),
    @("E44", "9.18%"),
    @("G44", "15"),
    @("D45", "0.01052"),
    @("E45", "32.23%"),
    @("G45", "15"),
    @("D46", "0.00006752"),
    @("E46", "2.77%"),
    @("G46", "15"),
    @("D47", "0.00000000752"),
    @("E47", "0.21%"),
    @("G47", "15"),
    @("D48", "0.002878"),
    @("E48", "-12.61%"),
    @("G48", "15"),
    @("E49", "59.31%"),
    @("G49", "15"),
    @("D50", "0.00002105"),
    @("E50", "0.21%"),
    @("G50", "15"),
    @("D51", "0.0002005"),
    @("E51", "0.21%"),
    @("G51", "15")
)

foreach ($edit in $edits) {
    Set-TextValue $ws $edit[0] $edit[1]
}

Write-Output "Applied $($edits.Count) cell updates"